$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 (pushing existing rows 20-43 down to 21-44).
# First copy row 20 so the newly inserted row inherits full formatting/values,
# then insert it, which shifts the original row 20 (and everything below) down by one.
$ws.Rows.Item(20).Copy()
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the new data values from the diff.
$ws.Cells.Item(20, 4).Value = 44810   # D20 - Fecha
$ws.Cells.Item(20, 10).Value = 300    # J20 - Volumen
$ws.Cells.Item(20, 11).Value = 2000   # K20 - Precio minimo
$ws.Cells.Item(20, 12).Value = 2500   # L20 - Precio maximo
$ws.Cells.Item(20, 13).Value = 2250   # M20 - Precio promedio ponderado
$ws.Cells.Item(20, 16).Value = 2250   # P20 - Precio $/Kg
